# Update crypto price/volume data as scraped by the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices formatted as plain text (values like "315.30" or
# "29.261.83" are not valid numbers/would lose trailing zeros if parsed),
# so force the column to Text before writing the new figures.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.261.83"
$ws.Range("E2").Value = "  +2.64%  "
$ws.Range("D3").Value = "1.894.25"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  -0.88%  "
$ws.Range("D5").Value = "315.30"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  -0.91%  "
$ws.Range("D7").Value = "0.5142"
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("D8").Value = "0.3911"
$ws.Range("E8").Value = "  -0.99%  "
$ws.Range("D9").Value = "0.08415"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").Value = "42.41"
$ws.Range("E10").Value = "  +1.43%  "
$ws.Range("D11").Value = "1.113"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.898.21"
$ws.Range("E12").Value = "  +0.62%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "6.235"
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("D14").Value = "20.63"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").Value = "7.304"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "1.006"
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("D17").Value = "93.03"
$ws.Range("E17").Value = "  +1.90%  "
$ws.Range("D18").Value = "0.00001105"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("D19").Value = "0.06742"
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("E21").Value = "  -0.90%  "
$ws.Range("D22").Value = "6.009"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("D23").Value = "29.269.10"
$ws.Range("E23").Value = "  +2.54%  "
$ws.Range("D24").Value = "11.12"
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("D25").Value = "2.214"
$ws.Range("E25").Value = "  -2.71%  "
$ws.Range("D26").Value = "2.112.46"
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("D27").Value = "158.92"
$ws.Range("E27").Value = "  -1.47%  "
$ws.Range("D28").Value = "20.86"
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("D29").Value = "2.431"
$ws.Range("E29").Value = "  +1.87%  "
$ws.Range("D30").Value = "127.78"
$ws.Range("E30").Value = "  +0.90%  "
$ws.Range("D31").Value = "1.056"
$ws.Range("E31").Value = "  +0.61%  "
$ws.Range("D32").Value = "0.1045"
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("D33").Value = "6.120"
$ws.Range("E33").Value = "  +5.68%  "
$ws.Range("D34").Value = "3.654"
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("D35").Value = "0.02474"
$ws.Range("E35").Value = "  +1.35%  "
$ws.Range("D36").Value = "0.06527"
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("D37").Value = "9.037"
$ws.Range("E37").Value = "  +0.94%  "
$ws.Range("D38").Value = "0.2189"
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("D39").Value = "1.224"
$ws.Range("E39").Value = "  +2.45%  "
$ws.Range("D40").Value = "5.130"
$ws.Range("E40").Value = "  +0.77%  "
$ws.Range("D41").Value = "0.6492"
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("D42").Value = "1.231"
$ws.Range("E42").Value = "  -2.96%  "
$ws.Range("D43").Value = "11.25"
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("D44").Value = "0.6049"
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("D45").Value = "13.20"
$ws.Range("E45").Value = "  +0.85%  "
$ws.Range("D46").Value = "3.679"
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("D47").Value = "2.041"
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("D49").Value = "122.83"
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("D50").Value = "1.178"
$ws.Range("E50").Value = "  -2.23%  "
$ws.Range("D51").Value = "77.46"
$ws.Range("E51").Value = "  +0.38%  "
